# Sample Project / Main.xlsx edit:
# Cell B11 (sheet "Main") held the shared string "R40"; it is changed to the
# (text) value "1". A leading apostrophe forces Excel to store it as text
# instead of coercing it to the number 1, which adds a new shared-string
# table entry - matching the new <si><t>1</t></si> entry added to
# xl/sharedStrings.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
